$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 (DATE_TYPE_CODE): "004" -> "003"
# Leading apostrophe forces Excel to keep the leading-zero numeric-looking
# string as literal text instead of coercing it to the number 3.
$ws.Range("J2").Value = "'003"

# N2 (REPORT_DATE) text timestamp
$ws.Range("N2").Value = "2020-03-31 00:00:00"

# O2..V2 numeric metrics
$ws.Range("O2").Value = 24299416.77
$ws.Range("P2").Value = 129.6794335463
$ws.Range("Q2").Value = 49170659.41
$ws.Range("R2").Value = 262.4105475347
$ws.Range("S2").Value = 6623688.06
$ws.Range("T2").Value = 35.3488367124
$ws.Range("U2").Value = -5561643.87
$ws.Range("V2").Value = -29.6809933125

# Y2, Z2 numeric metrics
$ws.Range("Y2").Value = 1197973.05
$ws.Range("Z2").Value = 6.3932590645

# AA2, AB2 (NETCASH_FINANCE / NETCASH_FINANCE_RATIO) now blank.
# A lone leading apostrophe yields an explicit empty text value (matching
# the other blank cells in this row) rather than deleting the cell outright.
$ws.Range("AA2").Value = "'"
$ws.Range("AB2").Value = "'"

# AC2, AD2 numeric metrics
$ws.Range("AC2").Value = 18738065.17
$ws.Range("AD2").Value = 531.2287653855
